$wb = $excel.ActiveWorkbook

# ALC row 121: Mindful Medicine | Tincture of Mind
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1282
$ws.Range("J121").Value = 1387.3334
$ws.Range("L121").Value = 4162.0002
$ws.Range("N121").Value = -7656.0002

# ALC row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1617.7576
$ws.Range("I132").Value = 1148.4445
$ws.Range("J132").Value = 3729.6667
$ws.Range("K132").Value = 3445.3335
$ws.Range("L132").Value = 11189.0001
$ws.Range("M132").Value = -915.3335000000002
$ws.Range("N132").Value = -16249.0001

# ALC row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2900.913
$ws.Range("I138").Value = 1984.8462
$ws.Range("J138").Value = 4091.8
$ws.Range("K138").Value = 5954.5386
$ws.Range("L138").Value = 12275.4
$ws.Range("M138").Value = -814.5385999999999
$ws.Range("N138").Value = -22555.4

# ARM row 32: Ingot We Trust | Steel Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28917.576
$ws.Range("I32").Value = 5514.636
$ws.Range("K32").Value = 5514.636
$ws.Range("M32").Value = -5227.636

# ARM row 82: Belle of the Brawl | Titanium Vambraces of Fending
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H82").Value = 27681
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 27681
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 27681
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -28403

# ARM row 85: Shouldering the Shut-ins (L) | Titanium Vambraces of Fending
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H85").Value = 27681
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 27681
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 27681
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -30177

# ARM row 122: Haste for High Durium | High Durium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1493.2094
$ws.Range("I122").Value = 1431.9025
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 4295.7075
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -1845.7075
$ws.Range("N122").Value = -13150

# ARM row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 101929.5
$ws.Range("I132").Value = 148792.12
$ws.Range("J132").Value = 2346.4375
$ws.Range("K132").Value = 446376.36
$ws.Range("L132").Value = 7039.3125
$ws.Range("M132").Value = -443846.36
$ws.Range("N132").Value = -12099.3125

# BSM row 134: Ruthenium Supremium | Ruthenium Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 126695.875
$ws.Range("I134").Value = 159616.89
$ws.Range("J134").Value = 1596
$ws.Range("K134").Value = 478850.67
$ws.Range("L134").Value = 4788
$ws.Range("M134").Value = -476315.67
$ws.Range("N134").Value = -9858

# CRP row 31: Wall Not Found | Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2314.25
$ws.Range("I31").Value = 2015.5
$ws.Range("J31").Value = 2911.75
$ws.Range("K31").Value = 2015.5
$ws.Range("L31").Value = 2911.75
$ws.Range("M31").Value = -1720.5
$ws.Range("N31").Value = -3501.75

# CRP row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2314.25
$ws.Range("I34").Value = 2015.5
$ws.Range("J34").Value = 2911.75
$ws.Range("K34").Value = 2015.5
$ws.Range("L34").Value = 2911.75
$ws.Range("M34").Value = -1813.5
$ws.Range("N34").Value = -3315.75

# CRP row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 971.8276
$ws.Range("I58").Value = 1039.0526
$ws.Range("J58").Value = 844.1
$ws.Range("K58").Value = 1039.0526
$ws.Range("L58").Value = 844.1
$ws.Range("M58").Value = -836.0526
$ws.Range("N58").Value = -1250.1

# CRP row 105: Zelkova, My Love | Zelkova Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 4793.923
$ws.Range("I105").Value = 5840.1665
$ws.Range("J105").Value = 3897.1428
$ws.Range("K105").Value = 5840.1665
$ws.Range("L105").Value = 3897.1428
$ws.Range("M105").Value = -4093.1665
$ws.Range("N105").Value = -7391.1428

# CRP row 132: Hull Lotta Damage | Ginseng Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1880.3
$ws.Range("I132").Value = 1681.36
$ws.Range("J132").Value = 2875
$ws.Range("K132").Value = 5044.08
$ws.Range("L132").Value = 8625
$ws.Range("M132").Value = -2514.08
$ws.Range("N132").Value = -13685

# CRP row 134: Wood You Be Quiet | Ceiba Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4468.5674
$ws.Range("I134").Value = 4907.433
$ws.Range("J134").Value = 2587.7144
$ws.Range("K134").Value = 14722.299
$ws.Range("L134").Value = 7763.1432
$ws.Range("M134").Value = -12187.299
$ws.Range("N134").Value = -12833.1432

# CRP row 136: Turali Quality | Dark Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 971.8276
$ws.Range("I136").Value = 1039.0526
$ws.Range("J136").Value = 844.1
$ws.Range("K136").Value = 3117.1578
$ws.Range("L136").Value = 2532.3
$ws.Range("M136").Value = -567.1578
$ws.Range("N136").Value = -7632.3

# CUL row 50: Moving Up in the World | Rolanberry Cheese
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 329.2857
$ws.Range("I50").Value = 172.5
$ws.Range("J50").Value = 392
$ws.Range("K50").Value = 517.5
$ws.Range("L50").Value = 1176
$ws.Range("M50").Value = -36.5
$ws.Range("N50").Value = -2138

# CUL row 53: Rolanberry Fields Forever | Rolanberry Cheese
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 329.2857
$ws.Range("I53").Value = 172.5
$ws.Range("J53").Value = 392
$ws.Range("K53").Value = 517.5
$ws.Range("L53").Value = 1176
$ws.Range("M53").Value = -36.5
$ws.Range("N53").Value = -2138

# GSM row 102: Put the Metal to the Peddle | Durium Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1964.8
$ws.Range("I102").Value = 1705.6666
$ws.Range("J102").Value = 2353.5
$ws.Range("K102").Value = 1705.6666
$ws.Range("L102").Value = 2353.5
$ws.Range("M102").Value = -83.66660000000002
$ws.Range("N102").Value = -5597.5

# GSM row 132: On Board for Lar | Lar Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3660.6216
$ws.Range("I132").Value = 3447.4827
$ws.Range("K132").Value = 10342.4481
$ws.Range("M132").Value = -7812.4481

# LTW row 40: Best Served Toad | Toad Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2298.3333
$ws.Range("I40").Value = 1700
$ws.Range("K40").Value = 1700
$ws.Range("M40").Value = -1564

# LTW row 122: Hell on Leather | Gaja Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3270.6667
$ws.Range("I122").Value = 2380
$ws.Range("J122").Value = 3716
$ws.Range("K122").Value = 7140
$ws.Range("L122").Value = 11148
$ws.Range("M122").Value = -4690
$ws.Range("N122").Value = -16048

# LTW row 132: Tenets of Tanning | Silver Lobo Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1935.1143
$ws.Range("I132").Value = 1561.08
$ws.Range("J132").Value = 2870.2
$ws.Range("K132").Value = 4683.24
$ws.Range("L132").Value = 8610.599999999999
$ws.Range("M132").Value = -2153.24
$ws.Range("N132").Value = -13670.6

# WVR row 119: A Job Well Done | Dwarven Cotton Gaskins of Fending
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 25000
$ws.Range("J119").Value = 25000
$ws.Range("L119").Value = 25000
$ws.Range("N119").Value = -34676

# WVR row 122: Heavy Armoire | Dark Hempen Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4864.643
$ws.Range("I122").Value = 6136.143
$ws.Range("J122").Value = 3593.1428
$ws.Range("K122").Value = 18408.429
$ws.Range("L122").Value = 10779.4284
$ws.Range("M122").Value = -15958.429
$ws.Range("N122").Value = -15679.4284

# WVR row 132: Comfy Cabins | Snow Cotton Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2619.875
$ws.Range("I132").Value = 2334.75
$ws.Range("J132").Value = 3475.25
$ws.Range("K132").Value = 7004.25
$ws.Range("L132").Value = 10425.75
$ws.Range("M132").Value = -4474.25
$ws.Range("N132").Value = -15485.75

# WVR row 136: Weaving the Envelope | Sarcenet Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1768.5454
$ws.Range("I136").Value = 1570.7858
$ws.Range("J136").Value = 2876
$ws.Range("K136").Value = 4712.357400000001
$ws.Range("L136").Value = 8628
$ws.Range("M136").Value = -2162.357400000001
$ws.Range("N136").Value = -13728
